# Scheduled-runner market data refresh: updates currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H:N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets to the latest scraped values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2772.762
$ws.Range("J98").Value = 6665
$ws.Range("L98").Value = 6665
$ws.Range("N98").Value = -9661

$ws.Range("H99").Value = 142857420
$ws.Range("I99").Value = 407
$ws.Range("K99").Value = 1221
$ws.Range("M99").Value = 277

$ws.Range("H122").Value = 2772.762
$ws.Range("J122").Value = 6665
$ws.Range("L122").Value = 19995
$ws.Range("N122").Value = -24895

$ws.Range("H132").Value = 1998.5405
$ws.Range("I132").Value = 1908.5518
$ws.Range("K132").Value = 5725.6554
$ws.Range("M132").Value = -3195.6554

$ws.Range("H137").Value = 2180210.8
$ws.Range("I137").Value = 3960
$ws.Range("K137").Value = 11880
$ws.Range("M137").Value = -9330

$ws.Range("H141").Value = 1693.45
$ws.Range("I141").Value = 1677.3158
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 5031.9474
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 148.0526
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 100000
$ws.Range("J24").Value = 100000
$ws.Range("L24").Value = 100000
$ws.Range("N24").Value = -100748

$ws.Range("H45").Value = 73594.21000000001
$ws.Range("I45").Value = 85359.914
$ws.Range("K45").Value = 85359.914
$ws.Range("M45").Value = -84982.914

$ws.Range("H46").Value = 3249.5
$ws.Range("I46").Value = 1999
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 1999
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -1680
$ws.Range("N46").Value = -5138

$ws.Range("H61").Value = 794757.9399999999
$ws.Range("I61").Value = 24083.725
$ws.Range("J61").Value = 2306465
$ws.Range("K61").Value = 24083.725
$ws.Range("L61").Value = 2306465
$ws.Range("M61").Value = -23871.725
$ws.Range("N61").Value = -2306889

$ws.Range("H74").Value = 430751.62
$ws.Range("I74").Value = 2707.9395
$ws.Range("J74").Value = 1174196
$ws.Range("K74").Value = 2707.9395
$ws.Range("L74").Value = 1174196
$ws.Range("M74").Value = -1833.9395
$ws.Range("N74").Value = -1175944

$ws.Range("H77").Value = 430751.62
$ws.Range("I77").Value = 2707.9395
$ws.Range("J77").Value = 1174196
$ws.Range("K77").Value = 13539.6975
$ws.Range("L77").Value = 5870980
$ws.Range("M77").Value = -9171.6975
$ws.Range("N77").Value = -5879716

$ws.Range("H97").Value = 7414.1113
$ws.Range("I97").Value = 8617.143
$ws.Range("K97").Value = 8617.143
$ws.Range("M97").Value = -8121.143

$ws.Range("H100").Value = 100000
$ws.Range("J100").Value = 100000
$ws.Range("L100").Value = 100000
$ws.Range("N100").Value = -102164

$ws.Range("H105").Value = 91687.22
$ws.Range("J105").Value = 91185
$ws.Range("L105").Value = 91185
$ws.Range("N105").Value = -98173

$ws.Range("H122").Value = 3179.625
$ws.Range("I122").Value = 3205.2856
$ws.Range("K122").Value = 9615.856800000001
$ws.Range("M122").Value = -7165.856800000001

$ws.Range("H136").Value = 794757.9399999999
$ws.Range("I136").Value = 24083.725
$ws.Range("J136").Value = 2306465
$ws.Range("K136").Value = 72251.17499999999
$ws.Range("L136").Value = 6919395
$ws.Range("M136").Value = -69701.17499999999
$ws.Range("N136").Value = -6924495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 28126778
$ws.Range("I134").Value = 1471.9546
$ws.Range("K134").Value = 4415.8638
$ws.Range("M134").Value = -1880.8638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2467.8948
$ws.Range("I31").Value = 2579.3157
$ws.Range("J31").Value = 2356.4736
$ws.Range("K31").Value = 2579.3157
$ws.Range("L31").Value = 2356.4736
$ws.Range("M31").Value = -2284.3157
$ws.Range("N31").Value = -2946.4736

$ws.Range("H34").Value = 2467.8948
$ws.Range("I34").Value = 2579.3157
$ws.Range("J34").Value = 2356.4736
$ws.Range("K34").Value = 2579.3157
$ws.Range("L34").Value = 2356.4736
$ws.Range("M34").Value = -2377.3157
$ws.Range("N34").Value = -2760.4736

$ws.Range("H58").Value = 1563.9615
$ws.Range("I58").Value = 1402.0714
$ws.Range("J58").Value = 1752.8334
$ws.Range("K58").Value = 1402.0714
$ws.Range("L58").Value = 1752.8334
$ws.Range("M58").Value = -1199.0714
$ws.Range("N58").Value = -2158.8334

$ws.Range("H132").Value = 2854.4
$ws.Range("I132").Value = 2748.75
$ws.Range("J132").Value = 2924.8333
$ws.Range("K132").Value = 8246.25
$ws.Range("L132").Value = 8774.499899999999
$ws.Range("M132").Value = -5716.25
$ws.Range("N132").Value = -13834.4999

$ws.Range("H134").Value = 2630.15
$ws.Range("I134").Value = 2474.6365
$ws.Range("J134").Value = 2820.2222
$ws.Range("K134").Value = 7423.9095
$ws.Range("L134").Value = 8460.6666
$ws.Range("M134").Value = -4888.9095
$ws.Range("N134").Value = -13530.6666

$ws.Range("H136").Value = 1563.9615
$ws.Range("I136").Value = 1402.0714
$ws.Range("J136").Value = 1752.8334
$ws.Range("K136").Value = 4206.2142
$ws.Range("L136").Value = 5258.5002
$ws.Range("M136").Value = -1656.2142
$ws.Range("N136").Value = -10358.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 9874451
$ws.Range("J9").Value = 12613633
$ws.Range("L9").Value = 37840899
$ws.Range("N9").Value = -37841347

$ws.Range("H98").Value = 595.6667
$ws.Range("J98").Value = 664.8
$ws.Range("L98").Value = 1994.4
$ws.Range("N98").Value = -4990.4

$ws.Range("H123").Value = 3999.1667
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H139").Value = 3587.9285
$ws.Range("I139").Value = 2291.7273
$ws.Range("K139").Value = 6875.1819
$ws.Range("M139").Value = -1735.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 34073.4
$ws.Range("J39").Value = 34073.4
$ws.Range("L39").Value = 34073.4
$ws.Range("N39").Value = -35137.4

$ws.Range("H93").Value = 124749.5
$ws.Range("J93").Value = 124749.5
$ws.Range("L93").Value = 124749.5
$ws.Range("N93").Value = -128493.5

$ws.Range("H102").Value = 1957.3636
$ws.Range("I102").Value = 1917.75
$ws.Range("K102").Value = 1917.75
$ws.Range("M102").Value = -295.75

$ws.Range("H132").Value = 1074635.2
$ws.Range("I132").Value = 19553.834
$ws.Range("K132").Value = 58661.50199999999
$ws.Range("M132").Value = -56131.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1194.4286
$ws.Range("I16").Value = 1194.4286
$ws.Range("K16").Value = 1194.4286
$ws.Range("M16").Value = -1024.4286

$ws.Range("H94").Value = 39723.75
$ws.Range("J94").Value = 39723.75
$ws.Range("L94").Value = 39723.75
$ws.Range("N94").Value = -41075.75

$ws.Range("H98").Value = 42765
$ws.Range("J98").Value = 42765
$ws.Range("L98").Value = 42765
$ws.Range("N98").Value = -48755

$ws.Range("H132").Value = 2905.8635
$ws.Range("I132").Value = 2927.9412
$ws.Range("K132").Value = 8783.8236
$ws.Range("M132").Value = -6253.8236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H100").Value = 1130.5385
$ws.Range("I100").Value = 609.7273
$ws.Range("K100").Value = 1219.4546
$ws.Range("M100").Value = -678.4546

$ws.Range("H107").Value = 3177035.2
$ws.Range("I107").Value = 1996.25
$ws.Range("K107").Value = 5988.75
$ws.Range("M107").Value = -4068.75

$ws.Range("H132").Value = 1729.7693
$ws.Range("I132").Value = 1284.1072
$ws.Range("J132").Value = 2864.182
$ws.Range("K132").Value = 3852.3216
$ws.Range("L132").Value = 8592.545999999998
$ws.Range("M132").Value = -1322.3216
$ws.Range("N132").Value = -13652.546
